$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from an existing header cell to the new F1 header cell
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Fecha_Scraping"

# Fill F2:F62 with the scraping timestamp as text
$fecha = "2025-05-11 19:22:43"
for ($r = 2; $r -le 62; $r++) {
    $ws.Cells.Item($r, 6).Value = $fecha
}
